$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.970.41"
$ws.Range("E2").Value = "  +2.79%  "
$ws.Range("D3").Value = "3.818.74"
$ws.Range("E3").Value = "  +1.27%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'700.68"
$ws.Range("E5").Value = "  +8.33%  "
$ws.Range("D6").Value = "'173.13"
$ws.Range("E6").Value = "  +4.41%  "
$ws.Range("D7").Value = "3.818.21"
$ws.Range("E7").Value = "  +1.31%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  +0.64%  "
$ws.Range("E10").Value = "  +2.13%  "
$ws.Range("E11").Value = "  +4.84%  "
$ws.Range("E12").Value = "  +0.81%  "
$ws.Range("D14").Value = "'36.26"
$ws.Range("E14").Value = "  +3.76%  "
$ws.Range("D15").Value = "4.464.10"
$ws.Range("E15").Value = "  +1.58%  "
$ws.Range("D16").Value = "3.831.73"
$ws.Range("E16").Value = "  +1.94%  "
$ws.Range("D17").Value = "70.962.30"
$ws.Range("E17").Value = "  +2.90%  "
$ws.Range("D18").Value = "'17.77"
$ws.Range("E18").Value = "  -0.15%  "
$ws.Range("D19").Value = "'7.20"
$ws.Range("E19").Value = "  +2.80%  "
$ws.Range("D21").Value = "'11.18"
$ws.Range("E21").Value = "  +16.67%  "
$ws.Range("D22").Value = "'479.55"
$ws.Range("E22").Value = "  +2.78%  "
$ws.Range("D23").Value = "'0.712"
$ws.Range("E23").Value = "  +0.40%  "
$ws.Range("D24").Value = "'83.86"
$ws.Range("E24").Value = "  +2.58%  "
$ws.Range("E25").Value = "  +0.40%  "
$ws.Range("D26").Value = "'12.37"
$ws.Range("E26").Value = "  +1.16%  "
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").Value = "'10.46"
$ws.Range("E27").Value = "  +2.93%  "
$ws.Range("B28").Value = "Fetch.AI"
$ws.Range("C28").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D28").Value = "'2.16"
$ws.Range("E28").Value = "  +3.26%  "
$ws.Range("D29").Value = "3.971.75"
$ws.Range("E29").Value = "  +1.59%  "
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("D31").Value = "'3.10"
$ws.Range("E31").Value = "  +15.33%  "
$ws.Range("D32").Value = "'2.30"
$ws.Range("E32").Value = "  +2.06%  "
$ws.Range("D33").Value = "'7.52"
$ws.Range("E33").Value = "  +5.35%  "
$ws.Range("E34").Value = "  +6.37%  "
$ws.Range("D35").Value = "'29.56"
$ws.Range("E35").Value = "  +3.35%  "
$ws.Range("D36").Value = "'9.25"
$ws.Range("E36").Value = "  +4.78%  "
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("E38").Value = "  +2.82%  "
$ws.Range("D39").Value = "'3.41"
$ws.Range("E39").Value = "  +5.82%  "
$ws.Range("D40").Value = "'6.03"
$ws.Range("E40").Value = "  +4.46%  "
$ws.Range("D41").Value = "'2.24"
$ws.Range("E41").Value = "  +13.73%  "
$ws.Range("D42").Value = "'0.983"
$ws.Range("E42").Value = "  +3.01%  "
$ws.Range("B43").Value = "FLOKI"
$ws.Range("C43").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D43").Value = "'0.000325"
$ws.Range("E43").Value = "  +21.52%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").Value = "'162.73"
$ws.Range("E46").Value = "  +4.74%  "
$ws.Range("D47").Value = "'48.93"
$ws.Range("E47").Value = "  +3.36%  "
$ws.Range("D48").Value = "'44.48"
$ws.Range("E48").Value = "  -1.36%  "
$ws.Range("B49").Value = "TheGraph"
$ws.Range("C49").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D49").Value = "'0.301"
$ws.Range("E49").Value = "  +1.92%  "
$ws.Range("B50").Value = "ONDO"
$ws.Range("C50").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D50").Value = "'1.38"
$ws.Range("E50").Value = "  -1.53%  "
$ws.Range("B51").Value = "Bittensor"
$ws.Range("C51").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D51").Value = "'410.36"
$ws.Range("E51").Value = "  +7.18%  "
